# pyFIFOtax example_rsu.xlsx fix: "APPL" should always have been "AAPL"
$wb = $excel.ActiveWorkbook

# Fix the mis-typed ticker symbol on the "rsu" sheet (row 6, column B)
$wsRsu = $wb.Worksheets.Item("rsu")
$wsRsu.Range("B6").Value = "AAPL"

# Fix the same mis-typed ticker symbol on the "sell_orders" sheet (row 6, column B)
$wsSellOrders = $wb.Worksheets.Item("sell_orders")
$wsSellOrders.Range("B6").Value = "AAPL"

# Update the view state: rsu sheet now has B7 selected (but is not the active tab)
$wsRsu.Range("B7").Select()

# Move the active tab from "sell_orders" to "money_transfers"
$wsMoneyTransfers = $wb.Worksheets.Item("money_transfers")
$wsMoneyTransfers.Activate()
